$d = $word.ActiveDocument

foreach ($para in $d.Paragraphs) {
    $r = $para.Range
    if ($r.Text.TrimEnd([char]13, [char]7) -eq "Project Description") {
        $end = $r.End - 1
        $insertRange = $d.Range($end, $end)
        $insertRange.InsertAfter(" ")
        $insertRange.Font.Bold = $true
        $insertRange.Font.Size = 15
        break
    }
}
